# Add two new instruction rows ("smaller than" / "larget than") to the
# instruction list on Sheet1, and move the sheet selection to the last
# new cell (B24) - matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: ID 21 -> "smaller than" / "c=a<b;"
$ws.Range("A23").Value = 21
$ws.Range("C23").Value = "c=a<b;"

# Row 24: ID 22 -> "larget than" / "c=a>b;"
$ws.Range("A24").Value = 22
$ws.Range("C24").Value = "c=a>b;"

# Column B filled in after C so the shared-string table gets the same
# insertion order as the source commit (c=a<b;, c=a>b;, smaller than, larget than).
$ws.Range("B23").Value = "smaller than"
$ws.Range("B24").Value = "larget than"

# Move the active selection to the last cell touched.
$ws.Range("B24").Select()
